# Generate Report for Handoff
#
# Refreshes the localization-status report for the rows that are
# "Ready for handoff" (and don't yet have a handback file): rows
# 7, 8, 9, 10, 12, 14 on the zh-cn / de-de / Overview sheets.
#   - Overview!G  (Latest HO Xliff Generate Date)      -> new generate timestamp
#   - zh-cn!H     (Latest Handoff Datetime)             -> new handoff timestamp
#   - de-de!H     (Latest Handoff Datetime)              -> new generate timestamp (shared w/ Overview!G)
#   - zh-cn!E and de-de!E (Priority)                    -> "ht"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 12, 14)

$generateDate = "2016-08-23 14:23:22"
$handoffDate = "2016-08-23 14:23:13"

foreach ($r in $rows) {
    $overview.Range("G$r").Value = $generateDate
    $dede.Range("H$r").Value = $generateDate

    $zhcn.Range("H$r").Value = $handoffDate

    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
